# LOB1229.xlsx: "Trabalho de Graduacao II" syllabus refresh (2022-01-09 build).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new "Docentes responsaveis" row (7455355 - Robson da Silva Rocha) is added
# right after the existing "4780627 - Ana Lucia Gabas Ferreira" row (row 13).
# Inserting a whole row there shifts every following row down by one and keeps
# each row's height/formatting (the engine carries both along automatically).
$ws.Rows.Item(14).Insert()

# Three rows (Objectives-EN, Short-syllabus-EN, Syllabus-EN) previously only had
# column A populated; they now also need matching B/C column formatting before
# their new English translations are written in below. Clone it from row 10,
# which already carries the correct B (style 2) / C (style 3) formatting.
$ws.Range("B10:C10").Copy()
$ws.Range("B11:C11").PasteSpecial(-4122)
$ws.Range("B16:C16").PasteSpecial(-4122)
$ws.Range("B18:C18").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Field value updates ---
# Row 6: Creditos-trabalho: 2 -> 1
$ws.Cells.Item(6,2).Value = "1"
$ws.Cells.Item(6,3).Value = "1"

# Row 7: Carga horaria: 75 h -> 45 h
$ws.Cells.Item(7,2).Value = "45 h"
$ws.Cells.Item(7,3).Value = "45 h"

# Row 8: Ativacao: 01/01/2012 -> 01/01/2022
$ws.Cells.Item(8,2).Value = "01/01/2022"
$ws.Cells.Item(8,3).Value = "01/01/2022"

# Row 10: Objetivos (PT) rewritten for TG2
$ws.Cells.Item(10,2).Value = "O Trabalho de Graduação 2 tem por objetivo a integração, o aprofundamento e aplicação dos conhecimentos adquiridos ao longo do curso, preparando e desenvolvendo a capacidade do aluno para a realização de atividades, que fazem parte do perfil de atuação profissional do engenheiro ambiental. O desenvolvimento do Trabalho de Graduação 2 deverá seguir o modelo escolhido no Trabalho de Graduação 1, e permitirá o uso de tecnologias digitais ou outras metodologias para desenvolvimento de conteúdo ou produto aplicável que utilize conteúdo da Engenharia Ambiental ou correlatas."
$ws.Cells.Item(10,3).Value = "O Trabalho de Graduação 2 tem por objetivo a integração, o aprofundamento e aplicação dos conhecimentos adquiridos ao longo do curso, preparando e desenvolvendo a capacidade do aluno para a realização de atividades, que fazem parte do perfil de atuação profissional do engenheiro ambiental. O desenvolvimento do Trabalho de Graduação 2 deverá seguir o modelo escolhido no Trabalho de Graduação 1, e permitirá o uso de tecnologias digitais ou outras metodologias para desenvolvimento de conteúdo ou produto aplicável que utilize conteúdo da Engenharia Ambiental ou correlatas."

# Row 11: Objectives (EN) - new translation added
$ws.Cells.Item(11,2).Value = "Graduation Work 2 aims to integrate, deepen and apply the knowledge acquired throughout the course, preparing and developing the student's ability to carry out activities, which are part of the professional performance profile of the environmental engineer. The development of Graduate Work 2 should follow the model chosen in Graduate Work 1, and will allow the use of digital technologies or other methodologies for the development of content or applicable product that uses Environmental Engineering or related content."
$ws.Cells.Item(11,3).Value = "Graduation Work 2 aims to integrate, deepen and apply the knowledge acquired throughout the course, preparing and developing the student's ability to carry out activities, which are part of the professional performance profile of the environmental engineer. The development of Graduate Work 2 should follow the model chosen in Graduate Work 1, and will allow the use of digital technologies or other methodologies for the development of content or applicable product that uses Environmental Engineering or related content."

# Row 14: Docentes responsaveis - new co-professor row
$ws.Cells.Item(14,2).Value = "7455355 - Robson da Silva Rocha"
$ws.Cells.Item(14,3).Value = "7455355 - Robson da Silva Rocha"

# Row 15: Programa resumido (PT) rewritten
$ws.Cells.Item(15,2).Value = "O aluno deverá desenvolver o Trabalho de Graduação 2 seguindo o modelo escolhido e já utilizado no Trabalho de Graduação 1. Diante da escolha do modelo, o aluno deverá desenvolver conteúdo científico ou produto aplicável, ambos relacionados a Engenharia Ambiental, e ao final, o Trabalho de Graduação deverá ser avaliado por banca avaliadora."
$ws.Cells.Item(15,3).Value = "O aluno deverá desenvolver o Trabalho de Graduação 2 seguindo o modelo escolhido e já utilizado no Trabalho de Graduação 1. Diante da escolha do modelo, o aluno deverá desenvolver conteúdo científico ou produto aplicável, ambos relacionados a Engenharia Ambiental, e ao final, o Trabalho de Graduação deverá ser avaliado por banca avaliadora."

# Row 16: Short syllabus (EN) - new translation added
$ws.Cells.Item(16,2).Value = "The student must develop the Graduate Work 2 following the model chosen and already used in the Graduate Work 1. When choosing the model, the student must develop scientific content or applicable product, both related to Environmental Engineering, and at the end, the Graduate Work must be evaluated by an evaluating panel."
$ws.Cells.Item(16,3).Value = "The student must develop the Graduate Work 2 following the model chosen and already used in the Graduate Work 1. When choosing the model, the student must develop scientific content or applicable product, both related to Environmental Engineering, and at the end, the Graduate Work must be evaluated by an evaluating panel."

# Row 17: Programa (PT) rewritten
$ws.Cells.Item(17,2).Value = "Para ambos os modelos (artigo ou produto): O aluno deverá dar continuidade ao desenvolvimento do Trabalho iniciado no Trabalho de Graduação 1. O programa da disciplina é constituído pelas seguintes etapas: 1) Desenvolvimento do tema com base nas atividades de cronograma aprovados no Trabalho de Graduação 1. 2) Desenvolvimento do texto final, conforme modelo fornecido pelos responsáveis da disciplina. 3) Entrega da versão final do texto, com aprovação do orientador e com a indicação da banca de avaliação. 4) Avaliação e atribuição de notas pela banca avaliadora, que pode ou não ser a mesma composta anteriormente no Trabalho de Graduação 1, a critério do orientador."
$ws.Cells.Item(17,3).Value = "Para ambos os modelos (artigo ou produto): O aluno deverá dar continuidade ao desenvolvimento do Trabalho iniciado no Trabalho de Graduação 1. O programa da disciplina é constituído pelas seguintes etapas: 1) Desenvolvimento do tema com base nas atividades de cronograma aprovados no Trabalho de Graduação 1. 2) Desenvolvimento do texto final, conforme modelo fornecido pelos responsáveis da disciplina. 3) Entrega da versão final do texto, com aprovação do orientador e com a indicação da banca de avaliação. 4) Avaliação e atribuição de notas pela banca avaliadora, que pode ou não ser a mesma composta anteriormente no Trabalho de Graduação 1, a critério do orientador."

# Row 18: Syllabus (EN) - new translation added
$ws.Cells.Item(18,2).Value = "For both models (article or product): The student must continue with the development of the Work initiated in the Graduate Work 1. The course program consists of the following steps: 1) Development of the theme based on the schedule activities approved in the Graduate Work 1. 2) Development of the final text, according to the model provided by those responsible for the discipline. 3) Delivery of the final version of the text, with the approval of the advisor and with the indication of the evaluation board. 4) Evaluation and grading by the examining board, which may or may not be the same previously composed in Graduate Work 1, at the discretion of the advisor."
$ws.Cells.Item(18,3).Value = "For both models (article or product): The student must continue with the development of the Work initiated in the Graduate Work 1. The course program consists of the following steps: 1) Development of the theme based on the schedule activities approved in the Graduate Work 1. 2) Development of the final text, according to the model provided by those responsible for the discipline. 3) Delivery of the final version of the text, with the approval of the advisor and with the indication of the evaluation board. 4) Evaluation and grading by the examining board, which may or may not be the same previously composed in Graduate Work 1, at the discretion of the advisor."

# Row 20: Metodo (PT) rewritten
$ws.Cells.Item(20,2).Value = "Para ambos os modelos (artigo ou produto): O aluno deverá apresentar um artigo científico seguindo modelo fornecido pelos responsáveis da disciplina. A versão final do texto deverá ser aprovada pelo orientador no ato da submissão e deverá ser apresentado para banca de avaliação composta por dois doutores ou especialistas na área do projeto. Alternativamente, pode-se apresentar o texto com formatação de outra revista desde que seja anexado o comprovante de submissão do artigo, nesse caso, o texto apresentado pode seguir as regras de formatação da revista escolhida. A revista escolhida deve ser indexada por algum sistema de base de dados (Web of Science, Scopus, Scielo, Cinahl, Medline, etc). Alternativamente, artigos aprovados em revistas da área de estudo, até a semana anterior a apresentação, desobrigam o aluno a apresentar o trabalho para a banca e nesse caso, a entrega do artigo deve ser acompanhada pela comprovação do aceite do trabalho. A revista escolhida deve ser indexada por algum sistema de base de dados (Web of Science, Scopus, Scielo, Cinahl, Medline, etc)"
$ws.Cells.Item(20,3).Value = "Para ambos os modelos (artigo ou produto): O aluno deverá apresentar um artigo científico seguindo modelo fornecido pelos responsáveis da disciplina. A versão final do texto deverá ser aprovada pelo orientador no ato da submissão e deverá ser apresentado para banca de avaliação composta por dois doutores ou especialistas na área do projeto. Alternativamente, pode-se apresentar o texto com formatação de outra revista desde que seja anexado o comprovante de submissão do artigo, nesse caso, o texto apresentado pode seguir as regras de formatação da revista escolhida. A revista escolhida deve ser indexada por algum sistema de base de dados (Web of Science, Scopus, Scielo, Cinahl, Medline, etc). Alternativamente, artigos aprovados em revistas da área de estudo, até a semana anterior a apresentação, desobrigam o aluno a apresentar o trabalho para a banca e nesse caso, a entrega do artigo deve ser acompanhada pela comprovação do aceite do trabalho. A revista escolhida deve ser indexada por algum sistema de base de dados (Web of Science, Scopus, Scielo, Cinahl, Medline, etc)"

# Row 21: Criterio (PT) rewritten
$ws.Cells.Item(21,2).Value = "Avaliação e emissão de parecer pela banca avaliadora e pelo orientador, com atribuição de nota única final.Fica sob responsabilidade do orientador a verificação de ocorrência de plágio utilizando software apropriado e avaliação em Comitê de Ética, quando exigido, via cadastro na Plataforma Brasil."
$ws.Cells.Item(21,3).Value = "Avaliação e emissão de parecer pela banca avaliadora e pelo orientador, com atribuição de nota única final.Fica sob responsabilidade do orientador a verificação de ocorrência de plágio utilizando software apropriado e avaliação em Comitê de Ética, quando exigido, via cadastro na Plataforma Brasil."

# Row 23: Bibliografia (PT) rewritten
$ws.Cells.Item(23,2).Value = "A ser definido no decorrer de cada projeto"
$ws.Cells.Item(23,3).Value = "A ser definido no decorrer de cada projeto"
